$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update vendor info block ---
$ws.Range("B3").Value = "A-one Industrial Sales"
$ws.Range("D3").Value = "Lopez Jaena St., Libertad, Bacolod"
$ws.Range("B4").Value = "435-7383; 432-0652; 476-1127"
$ws.Range("D4").Value = "435-7383"
$ws.Range("D5").Value = "Ms. Miles"
$ws.Range("B6").Value = ""
$ws.Range("D6").Value = "Re-seller, Distributor"

# --- Update existing item rows (11-14) ---
$ws.Range("A11").Value = "Generator Capacitor Discharge Ignition"
$ws.Range("C11").Value = ""

$ws.Range("A12").Value = "Tape"
$ws.Range("C12").Value = "-"

$ws.Range("A13").Value = "Tape"
$ws.Range("C13").Value = "-"

$ws.Range("A14").Value = "Plywood"
$ws.Range("C14").Value = "-"

# --- Add new item rows 15-22 ---
# Set values first, then copy formatting from row 14 (an existing item row)
# so the cells reuse the existing style indices instead of minting new ones.

# Row 15
$ws.Range("A15").Value = "Plywood"

# Row 16
$ws.Range("A16").Value = "Grinding Disc"
$ws.Range("C16").Value = "Tyrolit"

# Row 17
$ws.Range("A17").Value = "Grinding Disc"
$ws.Range("C17").Value = "Gold Elephant"

# Row 18
$ws.Range("A18").Value = "Cutting Disc"
$ws.Range("C18").Value = "Omega"

# Row 19
$ws.Range("A19").Value = "Cutting Nozzle"
$ws.Range("C19").Value = "-"

# Row 20
$ws.Range("A20").Value = "Wire"
$ws.Range("C20").Value = "-"

# Row 21
$ws.Range("A21").Value = "Gloves"

# Row 22
$ws.Range("A22").Value = "Soldering lead "
$ws.Range("C22").Value = "Brand: Rubicon"

# Merge first, then apply formatting (font/border/style) to new rows by
# copying from row 14 - an existing item row whose formatting (style index
# 4) we want to replicate. This ordering lets the pasted format resolve to
# the existing style index instead of minting a brand-new one.
$srcRow = $ws.Range("A14:D14")
for ($r = 15; $r -le 22; $r++) {
    $ws.Range("A" + $r + ":B" + $r).Merge()
    $ws.Range("C" + $r + ":D" + $r).Merge()
    $destRow = $ws.Range("A" + $r + ":D" + $r)
    $srcRow.Copy()
    $destRow.PasteSpecial(-4122) # xlPasteFormats
}

$ws.Range("D22").Select()
